$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4363.203762799113
$ws.Range("C3").Value = 4363.203762799113
$ws.Range("C4").Value = 4363.203762799113
$ws.Range("C5").Value = 4363.203762799113
$ws.Range("C6").Value = 4277.040679972994
$ws.Range("C7").Value = 4277.040679972994
$ws.Range("C8").Value = 4271.792090327821
$ws.Range("C9").Value = 4211.169281233559
$ws.Range("C10").Value = 4211.169281233559
$ws.Range("C11").Value = 4186.956042960988
$ws.Range("C12").Value = 4182.224352183007
